$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item(1)

$ws.Range("H28").Value = 3331.9167
$ws.Range("I28").Value = 1564
$ws.Range("J28").Value = 5099.8335
$ws.Range("K28").Value = 1564
$ws.Range("L28").Value = 5099.8335
$ws.Range("M28").Value = -1079
$ws.Range("N28").Value = -6069.8335
$ws.Range("H40").Value = 2802.7222
$ws.Range("I40").Value = 1848.1666
$ws.Range("K40").Value = 1848.1666
$ws.Range("M40").Value = -1673.1666
$ws.Range("H62").Value = 6212.75
$ws.Range("I62").Value = 4925.6665
$ws.Range("K62").Value = 4925.6665
$ws.Range("M62").Value = -4301.6665
$ws.Range("H65").Value = 6212.75
$ws.Range("I65").Value = 4925.6665
$ws.Range("K65").Value = 24628.3325
$ws.Range("M65").Value = -21508.3325
$ws.Range("H70").Value = 2337.842
$ws.Range("I70").Value = 2076.7
$ws.Range("J70").Value = 2628
$ws.Range("K70").Value = 6230.099999999999
$ws.Range("L70").Value = 7884
$ws.Range("M70").Value = -5960.099999999999
$ws.Range("N70").Value = -8424
$ws.Range("H73").Value = 2337.842
$ws.Range("I73").Value = 2076.7
$ws.Range("J73").Value = 2628
$ws.Range("K73").Value = 6230.099999999999
$ws.Range("L73").Value = 7884
$ws.Range("M73").Value = -5294.099999999999
$ws.Range("N73").Value = -9756
$ws.Range("H98").Value = 502.46155
$ws.Range("I98").Value = 477.875
$ws.Range("J98").Value = 797.5
$ws.Range("K98").Value = 477.875
$ws.Range("L98").Value = 797.5
$ws.Range("M98").Value = 1020.125
$ws.Range("N98").Value = -3793.5
$ws.Range("H122").Value = 502.46155
$ws.Range("I122").Value = 477.875
$ws.Range("J122").Value = 797.5
$ws.Range("K122").Value = 1433.625
$ws.Range("L122").Value = 2392.5
$ws.Range("M122").Value = 1016.375
$ws.Range("N122").Value = -7292.5
$ws.Range("H137").Value = 2083.7097
$ws.Range("I137").Value = 2024.875
$ws.Range("K137").Value = 6074.625
$ws.Range("M137").Value = -3524.625
$ws.Range("H138").Value = 6807629.5
$ws.Range("I138").Value = 2428.2
$ws.Range("J138").Value = 8552553
$ws.Range("K138").Value = 7284.599999999999
$ws.Range("L138").Value = 25657659
$ws.Range("M138").Value = -2144.599999999999
$ws.Range("N138").Value = -25667939
# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item(2)

$ws.Range("H32").Value = 15634937
$ws.Range("I32").Value = 24394826
$ws.Range("K32").Value = 24394826
$ws.Range("M32").Value = -24394539
$ws.Range("H61").Value = 35718704
$ws.Range("I61").Value = 52634320
$ws.Range("K61").Value = 52634320
$ws.Range("M61").Value = -52634108
$ws.Range("H64").Value = 100000
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 100000
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H74").Value = 21765122
$ws.Range("I74").Value = 29445748
$ws.Range("K74").Value = 29445748
$ws.Range("M74").Value = -29444874
$ws.Range("H77").Value = 21765122
$ws.Range("I77").Value = 29445748
$ws.Range("K77").Value = 147228740
$ws.Range("M77").Value = -147224372
$ws.Range("H102").Value = 68193.97
$ws.Range("I102").Value = 75178.07000000001
$ws.Range("K102").Value = 75178.07000000001
$ws.Range("M102").Value = -73556.07000000001
$ws.Range("H132").Value = 37047480
$ws.Range("I132").Value = 11848.046
$ws.Range("K132").Value = 35544.138
$ws.Range("M132").Value = -33014.138
$ws.Range("H136").Value = 35718704
$ws.Range("I136").Value = 52634320
$ws.Range("K136").Value = 157902960
$ws.Range("M136").Value = -157900410
# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item(3)

$ws.Range("H86").Value = 19359.422
$ws.Range("I86").Value = 12122.091
$ws.Range("J86").Value = 29310.75
$ws.Range("K86").Value = 12122.091
$ws.Range("L86").Value = 29310.75
$ws.Range("M86").Value = -10999.091
$ws.Range("N86").Value = -31556.75
$ws.Range("H89").Value = 19359.422
$ws.Range("I89").Value = 12122.091
$ws.Range("J89").Value = 29310.75
$ws.Range("K89").Value = 60610.455
$ws.Range("L89").Value = 146553.75
$ws.Range("M89").Value = -54994.455
$ws.Range("N89").Value = -157785.75
$ws.Range("H94").Value = 1782.6666
$ws.Range("I94").Value = 1516.3871
$ws.Range("J94").Value = 2533.0908
$ws.Range("K94").Value = 1516.3871
$ws.Range("L94").Value = 2533.0908
$ws.Range("M94").Value = -1065.3871
$ws.Range("N94").Value = -3435.0908
$ws.Range("H105").Value = 8929.154
$ws.Range("I105").Value = 11019.4
$ws.Range("J105").Value = 1961.6666
$ws.Range("K105").Value = 11019.4
$ws.Range("L105").Value = 1961.6666
$ws.Range("M105").Value = -9272.4
$ws.Range("N105").Value = -5455.6666
$ws.Range("H134").Value = 3653.1538
$ws.Range("I134").Value = 3694.8696
$ws.Range("J134").Value = 3333.3333
$ws.Range("K134").Value = 11084.6088
$ws.Range("L134").Value = 9999.999899999999
$ws.Range("M134").Value = -8549.6088
$ws.Range("N134").Value = -15069.9999
# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item(4)

$ws.Range("H62").Value = 4497.5713
$ws.Range("I62").Value = 3374.5
$ws.Range("K62").Value = 3374.5
$ws.Range("M62").Value = -2750.5
$ws.Range("H65").Value = 4497.5713
$ws.Range("I65").Value = 3374.5
$ws.Range("K65").Value = 16872.5
$ws.Range("M65").Value = -13752.5
$ws.Range("H94").Value = 1920.4
$ws.Range("I94").Value = 1787.25
$ws.Range("K94").Value = 1787.25
$ws.Range("M94").Value = -1336.25
# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item(5)

$ws.Range("H133").Value = 12947
$ws.Range("I133").Value = 7600.6
$ws.Range("K133").Value = 22801.8
$ws.Range("M133").Value = -17741.8
$ws.Range("H134").Value = 8591.65
$ws.Range("J134").Value = 19499.857
$ws.Range("L134").Value = 58499.571
$ws.Range("N134").Value = -68639.571
# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item(6)

$ws.Range("H70").Value = 204681.4
$ws.Range("I70").Value = 404363.2
$ws.Range("J70").Value = 4999.6
$ws.Range("K70").Value = 404363.2
$ws.Range("L70").Value = 4999.6
$ws.Range("M70").Value = -404093.2
$ws.Range("N70").Value = -5539.6
$ws.Range("H73").Value = 204681.4
$ws.Range("I73").Value = 404363.2
$ws.Range("J73").Value = 4999.6
$ws.Range("K73").Value = 404363.2
$ws.Range("L73").Value = 4999.6
$ws.Range("M73").Value = -403427.2
$ws.Range("N73").Value = -6871.6
$ws.Range("H97").Value = 2199.6843
$ws.Range("I97").Value = 584.1
$ws.Range("J97").Value = 3994.7778
$ws.Range("K97").Value = 584.1
$ws.Range("L97").Value = 3994.7778
$ws.Range("M97").Value = -88.10000000000002
$ws.Range("N97").Value = -4986.7778
$ws.Range("H132").Value = 3110.4443
$ws.Range("I132").Value = 3187.1765
$ws.Range("J132").Value = 1806
$ws.Range("K132").Value = 9561.529500000001
$ws.Range("L132").Value = 5418
$ws.Range("M132").Value = -7031.529500000001
$ws.Range("N132").Value = -10478
# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item(7)

$ws.Range("H25").Value = 37575
$ws.Range("J25").Value = 37575
$ws.Range("L25").Value = 37575
$ws.Range("N25").Value = -38035
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H122").Value = 4066.5557
$ws.Range("I122").Value = 3577
$ws.Range("K122").Value = 10731
$ws.Range("M122").Value = -8281
$ws.Range("H132").Value = 117648510
$ws.Range("I132").Value = 1481.1666
$ws.Range("K132").Value = 4443.4998
$ws.Range("M132").Value = -1913.4998
# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item(8)

$ws.Range("H74").Value = 6721.4443
$ws.Range("J74").Value = 6936.625
$ws.Range("L74").Value = 6936.625
$ws.Range("N74").Value = -8808.625
$ws.Range("H77").Value = 6721.4443
$ws.Range("J77").Value = 6936.625
$ws.Range("L77").Value = 20809.875
$ws.Range("N77").Value = -30169.875
$ws.Range("H107").Value = 411.375
$ws.Range("I107").Value = 428.69232
$ws.Range("J107").Value = 336.33334
$ws.Range("K107").Value = 1286.07696
$ws.Range("L107").Value = 1009.00002
$ws.Range("M107").Value = 633.9230400000001
$ws.Range("N107").Value = -4849.00002
$ws.Range("H113").Value = 433.23077
$ws.Range("I113").Value = 317.38095
$ws.Range("K113").Value = 952.14285
$ws.Range("M113").Value = 1217.85715
$ws.Range("H122").Value = 50055480
$ws.Range("I122").Value = 52689772
$ws.Range("J122").Value = 3998
$ws.Range("K122").Value = 158069316
$ws.Range("L122").Value = 11994
$ws.Range("M122").Value = -158066866
$ws.Range("N122").Value = -16894